# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1046
$ws1.Range("F5").Value  = 2831
$ws1.Range("F7").Value  = 236
$ws1.Range("F8").Value  = 22
$ws1.Range("F9").Value  = 126
$ws1.Range("F10").Value = 82
$ws1.Range("F11").Value = 107
$ws1.Range("F12").Value = 2667
$ws1.Range("F13").Value = 877

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1046
$ws4.Range("F6").Value  = 2831
$ws4.Range("F8").Value  = 236
$ws4.Range("F9").Value  = 22
$ws4.Range("F11").Value = 126
$ws4.Range("F12").Value = 82
$ws4.Range("F13").Value = 107
$ws4.Range("F14").Value = 2667
$ws4.Range("F15").Value = 877
